$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 7333.6875
$ws.Range("I74").Value = 5423.75
$ws.Range("J74").Value = 7970.3335
$ws.Range("K74").Value = 5423.75
$ws.Range("L74").Value = 7970.3335
$ws.Range("M74").Value = -4487.75
$ws.Range("N74").Value = -9842.333500000001
$ws.Range("H77").Value = 7333.6875
$ws.Range("I77").Value = 5423.75
$ws.Range("J77").Value = 7970.3335
$ws.Range("K77").Value = 27118.75
$ws.Range("L77").Value = 39851.6675
$ws.Range("M77").Value = -22438.75
$ws.Range("N77").Value = -49211.6675
$ws.Range("H86").Value = 5106.9165
$ws.Range("I86").Value = 4640.8335
$ws.Range("K86").Value = 4640.8335
$ws.Range("M86").Value = -3517.8335
$ws.Range("H89").Value = 5106.9165
$ws.Range("I89").Value = 4640.8335
$ws.Range("K89").Value = 23204.1675
$ws.Range("M89").Value = -17588.1675
$ws.Range("H92").Value = 2825.5715
$ws.Range("I92").Value = 2752.7727
$ws.Range("J92").Value = 3092.5
$ws.Range("K92").Value = 2752.7727
$ws.Range("L92").Value = 3092.5
$ws.Range("M92").Value = -1504.7727
$ws.Range("N92").Value = -5588.5
$ws.Range("H113").Value = 4649.875
$ws.Range("I113").Value = 4066.3333
$ws.Range("K113").Value = 4066.3333
$ws.Range("M113").Value = -812.3332999999998
$ws.Range("H137").Value = 3104.1785
$ws.Range("I137").Value = 2605.5715
$ws.Range("J137").Value = 4600
$ws.Range("K137").Value = 7816.7145
$ws.Range("L137").Value = 13800
$ws.Range("M137").Value = -5266.7145
$ws.Range("N137").Value = -18900
$ws.Range("H138").Value = 3299.1738
$ws.Range("J138").Value = 3617.7715
$ws.Range("L138").Value = 10853.3145
$ws.Range("N138").Value = -21133.3145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1911.7703
$ws.Range("I32").Value = 1188.0149
$ws.Range("K32").Value = 1188.0149
$ws.Range("M32").Value = -901.0148999999999
$ws.Range("H61").Value = 2225.2334
$ws.Range("I61").Value = 1939.9166
$ws.Range("K61").Value = 1939.9166
$ws.Range("M61").Value = -1727.9166
$ws.Range("H74").Value = 2508.56
$ws.Range("I74").Value = 1499.6
$ws.Range("K74").Value = 1499.6
$ws.Range("M74").Value = -625.5999999999999
$ws.Range("H77").Value = 2508.56
$ws.Range("I77").Value = 1499.6
$ws.Range("K77").Value = 7498
$ws.Range("M77").Value = -3130
$ws.Range("H97").Value = 1041.25
$ws.Range("I97").Value = 962.381
$ws.Range("J97").Value = 1593.3334
$ws.Range("K97").Value = 962.381
$ws.Range("L97").Value = 1593.3334
$ws.Range("M97").Value = -466.381
$ws.Range("N97").Value = -2585.3334
$ws.Range("H102").Value = 4034.5
$ws.Range("I102").Value = 3769.3333
$ws.Range("J102").Value = 4564.8335
$ws.Range("K102").Value = 3769.3333
$ws.Range("L102").Value = 4564.8335
$ws.Range("M102").Value = -2147.3333
$ws.Range("N102").Value = -7808.8335
$ws.Range("H132").Value = 2217.1143
$ws.Range("I132").Value = 1846.1786
$ws.Range("K132").Value = 5538.5358
$ws.Range("M132").Value = -3008.5358
$ws.Range("H136").Value = 2225.2334
$ws.Range("I136").Value = 1939.9166
$ws.Range("K136").Value = 5819.7498
$ws.Range("M136").Value = -3269.7498

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1215.1111
$ws.Range("I86").Value = 1106.9
$ws.Range("J86").Value = 1350.375
$ws.Range("K86").Value = 1106.9
$ws.Range("L86").Value = 1350.375
$ws.Range("M86").Value = 16.09999999999991
$ws.Range("N86").Value = -3596.375
$ws.Range("H89").Value = 1215.1111
$ws.Range("I89").Value = 1106.9
$ws.Range("J89").Value = 1350.375
$ws.Range("K89").Value = 5534.5
$ws.Range("L89").Value = 6751.875
$ws.Range("M89").Value = 81.5
$ws.Range("N89").Value = -17983.875
$ws.Range("H94").Value = 1194.5834
$ws.Range("I94").Value = 1236.5714
$ws.Range("K94").Value = 1236.5714
$ws.Range("M94").Value = -785.5714
$ws.Range("H99").Value = 30694.652
$ws.Range("I99").Value = 34321.59
$ws.Range("J99").Value = 20418.334
$ws.Range("K99").Value = 34321.59
$ws.Range("L99").Value = 20418.334
$ws.Range("M99").Value = -32823.59
$ws.Range("N99").Value = -23414.334
$ws.Range("H105").Value = 1237.2106
$ws.Range("I105").Value = 1249
$ws.Range("J105").Value = 1204.2
$ws.Range("K105").Value = 1249
$ws.Range("L105").Value = 1204.2
$ws.Range("M105").Value = 498
$ws.Range("N105").Value = -4698.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9957.25
$ws.Range("I31").Value = 11035.5
$ws.Range("J31").Value = 8879
$ws.Range("K31").Value = 11035.5
$ws.Range("L31").Value = 8879
$ws.Range("M31").Value = -10740.5
$ws.Range("N31").Value = -9469
$ws.Range("H34").Value = 9957.25
$ws.Range("I34").Value = 11035.5
$ws.Range("J34").Value = 8879
$ws.Range("K34").Value = 11035.5
$ws.Range("L34").Value = 8879
$ws.Range("M34").Value = -10833.5
$ws.Range("N34").Value = -9283
$ws.Range("H62").Value = 2427.5715
$ws.Range("I62").Value = 2557.4
$ws.Range("J62").Value = 2103
$ws.Range("K62").Value = 2557.4
$ws.Range("L62").Value = 2103
$ws.Range("M62").Value = -1933.4
$ws.Range("N62").Value = -3351
$ws.Range("H65").Value = 2427.5715
$ws.Range("I65").Value = 2557.4
$ws.Range("J65").Value = 2103
$ws.Range("K65").Value = 12787
$ws.Range("L65").Value = 10515
$ws.Range("M65").Value = -9667
$ws.Range("N65").Value = -16755
$ws.Range("H134").Value = 9728.951999999999
$ws.Range("I134").Value = 10812.206
$ws.Range("J134").Value = 5125.125
$ws.Range("K134").Value = 32436.618
$ws.Range("L134").Value = 15375.375
$ws.Range("M134").Value = -29901.618
$ws.Range("N134").Value = -20445.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 335669
$ws.Range("I80").Value = 386796.62
$ws.Range("J80").Value = 202737.2
$ws.Range("K80").Value = 386796.62
$ws.Range("L80").Value = 202737.2
$ws.Range("M80").Value = -385798.62
$ws.Range("N80").Value = -204733.2
$ws.Range("H83").Value = 335669
$ws.Range("I83").Value = 386796.62
$ws.Range("J83").Value = 202737.2
$ws.Range("K83").Value = 1933983.1
$ws.Range("L83").Value = 1013686
$ws.Range("M83").Value = -1928991.1
$ws.Range("N83").Value = -1023670
$ws.Range("H97").Value = 1395.4
$ws.Range("I97").Value = 2514.5
$ws.Range("K97").Value = 2514.5
$ws.Range("M97").Value = -2018.5
$ws.Range("H102").Value = 2029.0488
$ws.Range("I102").Value = 2029.0488
$ws.Range("K102").Value = 2029.0488
$ws.Range("M102").Value = -407.0488
$ws.Range("H107").Value = 2688.875
$ws.Range("I107").Value = 3404.6667
$ws.Range("J107").Value = 2259.4
$ws.Range("K107").Value = 3404.6667
$ws.Range("L107").Value = 2259.4
$ws.Range("M107").Value = -1484.6667
$ws.Range("N107").Value = -6099.4
$ws.Range("H132").Value = 29423086
$ws.Range("I132").Value = 38469830
$ws.Range("K132").Value = 115409490
$ws.Range("M132").Value = -115406960

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 943
$ws.Range("J22").Value = 1116.6666
$ws.Range("L22").Value = 1116.6666
$ws.Range("N22").Value = -1706.6666
$ws.Range("H27").Value = 943
$ws.Range("J27").Value = 1116.6666
$ws.Range("L27").Value = 1116.6666
$ws.Range("N27").Value = -1330.6666
$ws.Range("H82").Value = 447.5
$ws.Range("I82").Value = 447.5
$ws.Range("K82").Value = 447.5
$ws.Range("M82").Value = -86.5
$ws.Range("H85").Value = 447.5
$ws.Range("I85").Value = 447.5
$ws.Range("K85").Value = 447.5
$ws.Range("M85").Value = 800.5
$ws.Range("H93").Value = 2668.875
$ws.Range("I93").Value = 2433.25
$ws.Range("J93").Value = 3375.75
$ws.Range("K93").Value = 2433.25
$ws.Range("L93").Value = 3375.75
$ws.Range("M93").Value = -1185.25
$ws.Range("N93").Value = -5871.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1670
$ws.Range("I100").Value = 1765.7142
$ws.Range("K100").Value = 3531.4284
$ws.Range("M100").Value = -2990.4284
$ws.Range("H107").Value = 2237.5
$ws.Range("I107").Value = 1227.16
$ws.Range("K107").Value = 3681.48
$ws.Range("M107").Value = -1761.48
$ws.Range("H122").Value = 1258.3667
$ws.Range("I122").Value = 1202.92
$ws.Range("K122").Value = 3608.76
$ws.Range("M122").Value = -1158.76
$ws.Range("H132").Value = 3692.8823
$ws.Range("I132").Value = 2707.7778
$ws.Range("K132").Value = 8123.3334
$ws.Range("M132").Value = -5593.3334
